# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Pomelo"
# at row 213, pushing the existing rows 213-231 down to 214-232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 213 (existing data shifts down).
$ws.Rows.Item(213).Insert()

# Fill in the new row 213 with the new weekly record.
$ws.Cells.Item(213, 1).Value  = 4
$ws.Cells.Item(213, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(213, 3).Value  = "Los Lagos"
$ws.Cells.Item(213, 4).Value  = 44610
$ws.Cells.Item(213, 5).Value  = 10
$ws.Cells.Item(213, 6).Value  = "Fruta"
$ws.Cells.Item(213, 7).Value  = 100102
$ws.Cells.Item(213, 8).Value  = "Cítricos"
$ws.Cells.Item(213, 9).Value  = 100102006
$ws.Cells.Item(213, 10).Value = "Pomelo"
$ws.Cells.Item(213, 11).Value = "Start Ruby"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 160
$ws.Cells.Item(213, 14).Value = 14000
$ws.Cells.Item(213, 15).Value = 14000
$ws.Cells.Item(213, 16).Value = 14000
$ws.Cells.Item(213, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(213, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(213, 19).Value = 1000
$ws.Cells.Item(213, 20).Value = 14
